$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) and "全部类型" (sheet4) both contain the same
# event listing data that needs the "想去人数" (want-to-go count)
# column F updated for rows 2 and 3.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 237
    $ws.Range("F3").Value = 303
}
